$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B7").Value = 2002837.36
$ws.Range("C7").Value = -55.66425751146264
$ws.Range("D7").Value = 1922
$ws.Range("E7").Value = 1922
$ws.Range("F7").Value = 1042.058980228928
$ws.Range("G7").Value = 7.54070316418376
